$d = $word.ActiveDocument

# The sentence in question currently reads (as a single run):
#   " Using the Lena image, apply two different mathematical
#     algorithms/equations. Explain what happened in terms of numbers
#     and visualization."
# and is immediately preceded, in the same paragraph, by a separate
# run containing just the period that ends "Let's play Photoshop.".
#
# The target edit inserts " to it" right after "algorithms/equations"
# (turning "...equations. Explain..." into "...equations to it. Explain...")
# and leaves that sentence split across three runs sharing identical
# run formatting.
#
# A plain Find/Replace or InsertAfter on this paragraph causes the
# engine's run-coalescing pass to merge every adjacent run that shares
# identical formatting - which would also swallow the preceding
# "." run into the edited text. To keep that "." run intact we give it
# a throwaway, reverted formatting nudge first so it is not considered
# for merging with its neighbour once the real edit happens.

# --- Step 1: insulate the preceding "." run from the coalescing pass ---
$dotFind = $d.Content
$dotFind.Find.Execute("Photoshop.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$dotStart = $dotFind.End - 1
$dotRange = $d.Range($dotStart, $dotStart + 1)
$dotRange.Font.Bold = $true

# --- Step 2: insert " to it" right after "algorithms/equations" ---
$targetFind = $d.Content
$targetFind.Find.Execute("algorithms/equations", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$insertPos = $targetFind.End
$insertRange = $d.Range($insertPos, $insertPos)
$insertRange.InsertAfter(" to it")

# --- Step 3: restore the "." run's real formatting (stays its own run) ---
$dotRange2 = $d.Range($dotStart, $dotStart + 1)
$dotRange2.Font.Bold = $false

# --- Step 4: split the sentence run into the three target runs by
#     nudging (and reverting) formatting on the newly inserted " to it" ---
$toItRange = $d.Range($insertPos, $insertPos + 6)
$toItRange.Font.Bold = $true
$toItRange.Font.Bold = $false
